$d = $word.ActiveDocument

# Locate the three target paragraphs by their distinctive text, rather than
# trusting a fixed index, in case paragraph numbering differs slightly.
$rf5Idx = $null
$rnf3Idx = $null
$rnf5Idx = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -match "^RF5:") {
        $rf5Idx = $i
    } elseif ($t -match "^RNF3: La comunicaci") {
        $rnf3Idx = $i
    } elseif ($t -match "^RNF") {
        # handles "RNF" + "5" + ": El sistema..." split across runs,
        # so just match the generic RNF5 paragraph about 80% uptime.
        if ($t -match "80% del tiempo") {
            $rnf5Idx = $i
        }
    }
}

# --- 1) RF5 paragraph: just add red font color, text/runs unchanged ---
# Use the full paragraph Range (including its end-of-paragraph mark) so the
# color lands on both the run(s) and the paragraph mark's rPr (w:pPr/w:rPr).
$p1 = $d.Paragraphs($rf5Idx)
$p1.Range.Font.Color = 255

# --- 2) RNF3 (HTTP) paragraph: just add red font color, single run ---
$p2 = $d.Paragraphs($rnf3Idx)
$p2.Range.Font.Color = 255

# --- 3) RNF5 paragraph: consolidate the 5 runs into a single run with the
#        full text, then apply red font color to text + paragraph mark ---
$p3 = $d.Paragraphs($rnf5Idx)
$r3 = $p3.Range
$rng3 = $d.Range($r3.Start, $r3.End - 1)
# Force Word to rebuild the run list by first swapping in different text,
# then writing the final text back as a single contiguous run.
$rng3.Text = "TEMP_PLACEHOLDER_XYZ"

$p3b = $d.Paragraphs($rnf5Idx)
$r3b = $p3b.Range
$rng3b = $d.Range($r3b.Start, $r3b.End - 1)
$rng3b.Text = "RNF5: El sistema debe estar disponible al menos el 80% del tiempo."

$p3c = $d.Paragraphs($rnf5Idx)
$r3c = $p3c.Range
$r3c.Font.Color = 255
